# The forecast table grew from 52 to 53 data rows. All of the
# y_0 / y_0_forecast / y_1 / y_1_forecast values were recomputed (bugfix to
# the forecast evaluation / simulated rt_data), and one additional earlier
# observation (date 39400) was prepended, shifting every later row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date-column formatting (the "s=2" style already used by every
# A-column cell in the table) down to the newly added A53 cell by copying
# the format from an existing date cell.
$ws.Range("A3").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# New data for rows 2..53:
#   col A = date_of_forecast (serial date)
#   col B = y_0 (year)
#   col C = y_0_forecast
#   col D = y_1 (year)
#   col E = y_1_forecast
$rows = @(
    @(39400, 2007, 5.896808312953783, 2008, 7.318442086255605),
    @(39583, 2008, 6.056254825277896, 2009, 4.566338461218011),
    @(39765, 2008, 7.441962824572235, 2009, 8.262942840582955),
    @(39948, 2009, 8.604123301398037, 2010, 8.260999835306727),
    @(40130, 2009, 6.277541464866987, 2010, 7.915558093865038),
    @(40310, 2010, 5.436647924209592, 2011, 4.862860110364875),
    @(40494, 2010, 6.535114773304773, 2011, 6.402044794134309),
    @(40676, 2011, 6.334380382529425, 2012, 5.829578861489648),
    @(40862, 2011, 5.12051970717502, 2012, 5.160655543917292),
    @(41044, 2012, 3.88993859232436, 2013, 3.4300351921007),
    @(41228, 2012, 3.65682115264816, 2013, 3.62442670409151),
    @(41409, 2013, 2.513767348245044, 2014, 2.479992751939486),
    @(41592, 2013, 2.943878639034381, 2014, 3.771815305047821),
    @(41774, 2014, 2.723916849952834, 2015, 2.551173534479334),
    @(41957, 2014, 1.172679597477866, 2015, 1.204704113773114),
    @(42137, 2015, 1.326505206336948, 2016, 1.211929054838756),
    @(42321, 2015, 2.961845079861303, 2016, 2.562199564969392),
    @(42503, 2016, 2.321003614014883, 2017, 2.610201636760778),
    @(42689, 2016, 2.508469427909898, 2017, 2.632055757778851),
    @(42867, 2017, 2.468891199411116, 2018, 2.593292206016984),
    @(43053, 2017, 3.523703831572056, 2018, 3.331068508781954),
    @(43145, 2018, 2.442951431721241, 2019, 2.874490643722805),
    @(43235, 2018, 3.133596157287766, 2019, 3.419422858788335),
    @(43326, 2018, 1.273091634877033, 2019, 1.947147525128035),
    @(43418, 2018, 1.178605266817589, 2019, 1.467147844249106),
    @(43510, 2019, 1.388614840712377, 2020, 1.965204165904111),
    @(43600, 2019, 0.08486825492834971, 2020, 0.9311475558545057),
    @(43691, 2019, 2.983246785467752, 2020, 2.844177684788551),
    @(43783, 2019, 3.047037961814492, 2020, 2.757054249287738),
    @(43875, 2020, 3.192446966735796, 2021, 2.739953120486138),
    @(43966, 2020, 2.405224065057476, 2021, 2.152035263856344),
    @(44068, 2020, -0.4891791466461126, 2021, 0.6596954352570572),
    @(44159, 2020, -0.2228847697281378, 2021, 2.031388171904314),
    @(44251, 2021, 2.00987424531256, 2022, 1.446144849547082),
    @(44341, 2021, -0.678826357714013, 2022, 0.02414656897629097),
    @(44432, 2021, -0.7941560676977599, 2022, 1.553668892952742),
    @(44525, 2021, -1.165854108406617, 2022, 2.209187654930855),
    @(44617, 2022, 3.122586570118835, 2023, 1.070586793467498),
    @(44706, 2022, 0.7010162698181555, 2023, -0.6730402944081559),
    @(44798, 2022, 2.34069710769782, 2023, 1.423433561452137),
    @(44890, 2022, 2.501311189006916, 2023, 2.819423640661167),
    @(44981, 2023, 2.639029809913129, 2024, 1.61912427117199),
    @(45071, 2023, 2.158153176293576, 2024, 1.298013848993262),
    @(45163, 2023, 0.9259311313598806, 2024, 1.321448002249337),
    @(45254, 2023, 0.6753076481029074, 2024, 0.6583865490665364),
    @(45345, 2024, 1.29262043481877, 2025, 1.577641784671369),
    @(45436, 2024, -0.1645072558042915, 2025, 0.1619933518385297),
    @(45534, 2024, 1.935025917091848, 2025, 1.814429727910216),
    @(45618, 2024, 2.039329803030121, 2025, 2.453497774623137),
    @(45713, 2025, 1.562469275846001, 2026, 1.276308453105246),
    @(45800, 2025, 3.131832690451031, 2026, 2.199676451050503),
    @(45891, 2025, 2.43119486791763, 2026, 2.795512766421537)
)

$arr = New-Object 'object[,]' $rows.Length, 5
for ($i = 0; $i -lt $rows.Length; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A2:E53").Value = $arr
